$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 updates
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 updates
$ws.Range("B12").Value = 189
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "179/252"
